# poblacion total,m y f por municipio
# Refresh the fuel-price listing (informe_madrid): new prices pulled in and
# the rows re-sorted ascending by precio_gasoleo_a, which reshuffles which
# station/address/locality triple lands on each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.497, "BALLENOIL", "AVENIDA FUENLABRADA, 6", "HUMANES DE MADRID"),
    @(1.514, "SHELL", "AVENIDA FUENLABRADA, 110", "HUMANES DE MADRID"),
    @(1.599, "BALLENOIL", "CALLE CABO RUFINO LAZARO, 7", "ROZAS DE MADRID (LAS)"),
    @(1.649, "CARREFOUR", "CARRETERA MADRID-LA CORUÑA KM. 22", "ROZAS DE MADRID (LAS)"),
    @(1.669, "GALP", "CTRA. N-VI km 21,700", "ROZAS DE MADRID (LAS)"),
    @(1.688, "REPSOL", "CALLE COPENHAGUES/N, S/N", "ROZAS DE MADRID (LAS)"),
    @(1.688, "REPSOL", "A-6 km 25,5", "ROZAS DE MADRID (LAS)"),
    @(1.688, "REPSOL", "CTRA. M-505 km 5,5", "ROZAS DE MADRID (LAS)"),
    @(1.688, "REPSOL", "CARRETERA M-505 km 5.5", "ROZAS DE MADRID (LAS)"),
    @(1.689, "CEPSA", "CARRETERA M-405 KM. 6", "HUMANES DE MADRID"),
    @(1.699, "REPSOL HUMANES", "AVENIDA LA INDUSTRIA, S/N", "HUMANES DE MADRID"),
    @(1.699, "CEPSA", "CARRETERA M-405 KM. 5,6", "HUMANES DE MADRID"),
    @(1.699, "REPSOL", "CL MADRID, 52", "HUMANES DE MADRID"),
    @(1.709, "BP VALDONAIRE", "CARRETERA AVD.DE LA INDUSTRIA KM. 15", "HUMANES DE MADRID"),
    @(1.709, "BP HUMANES - EL MOLINO", "AVENIDA DE LAS FLORES, 2", "HUMANES DE MADRID"),
    @(1.709, "REPSOL", "CARRETERA AVENIDA  DE LA INDUSTRIA , 46 KM. 1,1", "HUMANES DE MADRID"),
    @(1.709, "REPSOL", "CR A-6, 20,3", "ROZAS DE MADRID (LAS)"),
    @(1.709, "BP LAS ROZAS", "CL LAS CRUCES  S/N", "ROZAS DE MADRID (LAS)"),
    @(1.861, "COSTCO", "CALLE INNOVACIÓN, 19", "ROZAS DE MADRID (LAS)")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
